$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "New product" in column T (col 20), matching the formatting of S1
$ws.Range("T1").Value = "New product"
$ws.Range("S1").Copy()
$ws.Range("T1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Size column T to best-fit its header content (matches the other bestFit header columns)
$ws.Columns("T:T").ColumnWidth = 8.8

# Update the active selection to T7, as reflected in the saved view state
$ws.Range("T7").Select()
